$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.08569625208598
$ws.Range("C2").Value = 8.282214553375475
$ws.Range("D2").Value = 10.47992087735639
$ws.Range("F2").Value = 31.56644401273309
$ws.Range("G2").Value = 3.641049173929286
$ws.Range("J2").Value = 11.08646941670659
$ws.Range("M2").Value = 17.55973421367547
$ws.Range("O2").Value = 23.07794417772273

$ws.Range("B3").Value = 13.50127079956086
$ws.Range("C3").Value = 7.788697683597315
$ws.Range("D3").Value = 10.45788956088253
$ws.Range("F3").Value = 31.65945336606849
$ws.Range("G3").Value = 3.643319690848775
$ws.Range("J3").Value = 11.1206848908304
$ws.Range("M3").Value = 17.36119815080934
$ws.Range("O3").Value = 23.18253979538053

$ws.Range("B4").Value = 13.13014811736975
$ws.Range("C4").Value = 7.468026579051791
$ws.Range("D4").Value = 10.4457803461756
$ws.Range("F4").Value = 31.72637539553306
$ws.Range("G4").Value = 3.644787188391561
$ws.Range("J4").Value = 11.14336307037637
$ws.Range("M4").Value = 17.2404747215811
$ws.Range("O4").Value = 23.25358937053716

$ws.Range("B5").Value = 12.97603105341357
$ws.Range("C5").Value = 7.332951286009466
$ws.Range("D5").Value = 10.44120551168113
$ws.Range("F5").Value = 31.75610345621449
$ws.Range("G5").Value = 3.645403720494028
$ws.Range("J5").Value = 11.15302447756411
$ws.Range("M5").Value = 17.19162303804739
$ws.Range("O5").Value = 23.28425194025687

$ws.Range("B6").Value = 12.95027256418846
$ws.Range("C6").Value = 7.310257792362843
$ws.Range("D6").Value = 10.44046768101598
$ws.Range("F6").Value = 31.76118782615329
$ws.Range("G6").Value = 3.645507215207528
$ws.Range("J6").Value = 11.15465410559303
$ws.Range("M6").Value = 17.1835334198616
$ws.Range("O6").Value = 23.28944645091618

$ws.Range("B7").Value = 13.12808101506491
$ws.Range("C7").Value = 7.466222659181645
$ws.Range("D7").Value = 10.44571718765126
$ws.Range("F7").Value = 31.72676638624037
$ws.Range("G7").Value = 3.644795428112792
$ws.Range("J7").Value = 11.14349166748499
$ws.Range("M7").Value = 17.23981443481514
$ws.Range("O7").Value = 23.25399598591063

$ws.Range("B8").Value = 13.88685198306576
$ws.Range("C8").Value = 8.115733264298255
$ws.Range("D8").Value = 10.472032127019
$ws.Range("F8").Value = 31.59647037096072
$ws.Range("G8").Value = 3.641816852162388
$ws.Range("J8").Value = 11.09792035547159
$ws.Range("M8").Value = 17.4910631699878
$ws.Range("O8").Value = 23.11258722924893

$ws.Range("B9").Value = 15.26951509401634
$ws.Range("C9").Value = 9.24816590840835
$ws.Range("D9").Value = 10.53473894403968
$ws.Range("F9").Value = 31.41928750933944
$ws.Range("G9").Value = 3.63655544932558
$ws.Range("J9").Value = 11.02180400110742
$ws.Range("M9").Value = 17.99099766303315
$ws.Range("O9").Value = 22.88978800732234

$ws.Range("B10").Value = 16.21231519253035
$ws.Range("C10").Value = 9.992978513812725
$ws.Range("D10").Value = 10.58736604580336
$ws.Range("F10").Value = 31.33743698592567
$ws.Range("G10").Value = 3.633039369597429
$ws.Range("J10").Value = 10.97395844389041
$ws.Range("M10").Value = 18.35987523350441
$ws.Range("O10").Value = 22.75976576776934

$ws.Range("B11").Value = 16.62389726535774
$ws.Range("C11").Value = 10.31277431360853
$ws.Range("D11").Value = 10.6126810876371
$ws.Range("F11").Value = 31.31079918881649
$ws.Range("G11").Value = 3.631514878321143
$ws.Range("J11").Value = 10.95394628552947
$ws.Range("M11").Value = 18.52743263739229
$ws.Range("O11").Value = 22.70801707330245

$ws.Range("B12").Value = 16.77716603405245
$ws.Range("C12").Value = 10.43113311466425
$ws.Range("D12").Value = 10.62246028434073
$ws.Range("F12").Value = 31.30224289853202
$ws.Range("G12").Value = 3.630948314494822
$ws.Range("J12").Value = 10.94662033412968
$ws.Range("M12").Value = 18.59079764972409
$ws.Range("O12").Value = 22.68949236359691

$ws.Range("B13").Value = 16.74427342104767
$ws.Range("C13").Value = 10.40576438376547
$ws.Range("D13").Value = 10.62034565609497
$ws.Range("F13").Value = 31.30401746095625
$ws.Range("G13").Value = 3.631069857906664
$ws.Range("J13").Value = 10.94818688917881
$ws.Range("M13").Value = 18.5771554387555
$ws.Range("O13").Value = 22.69343422501205

$ws.Range("B14").Value = 16.63655917619117
$ws.Range("C14").Value = 10.32256671490045
$ws.Range("D14").Value = 10.61348178915799
$ws.Range("F14").Value = 31.31006453981811
$ws.Range("G14").Value = 3.631468052088167
$ws.Range("J14").Value = 10.95333851948098
$ws.Range("M14").Value = 18.53264770728976
$ws.Range("O14").Value = 22.70647152024316

$ws.Range("B15").Value = 16.57024128232448
$ws.Range("C15").Value = 10.27124868711644
$ws.Range("D15").Value = 10.60930245451967
$ws.Range("F15").Value = 31.31396810737398
$ws.Range("G15").Value = 3.631713352969538
$ws.Range("J15").Value = 10.95652689214705
$ws.Range("M15").Value = 18.50537286088404
$ws.Range("O15").Value = 22.71459698371355

$ws.Range("B16").Value = 16.18506023502262
$ws.Range("C16").Value = 9.971695769587434
$ws.Range("D16").Value = 10.58573891057096
$ws.Range("F16").Value = 31.33939168741814
$ws.Range("G16").Value = 3.633140502930387
$ws.Range("J16").Value = 10.9753015724387
$ws.Range("M16").Value = 18.34891571262274
$ws.Range("O16").Value = 22.76329719467869

$ws.Range("B17").Value = 15.94425619632424
$ws.Range("C17").Value = 9.783052494595006
$ws.Range("D17").Value = 10.57163223123172
$ws.Range("F17").Value = 31.35770761664794
$ws.Range("G17").Value = 3.634035180460321
$ws.Range("J17").Value = 10.98726831884579
$ws.Range("M17").Value = 18.25283532754106
$ws.Range("O17").Value = 22.79507363169778

$ws.Range("B18").Value = 15.80412844522346
$ws.Range("C18").Value = 9.67275939802879
$ws.Range("D18").Value = 10.56364813440228
$ws.Range("F18").Value = 31.36923947858092
$ws.Range("G18").Value = 3.634556837125801
$ws.Range("J18").Value = 10.99431626973008
$ws.Range("M18").Value = 18.19755212330399
$ws.Range("O18").Value = 22.81404660432276

$ws.Range("B19").Value = 15.75640804802702
$ws.Range("C19").Value = 9.635108756828629
$ws.Range("D19").Value = 10.56096727203555
$ws.Range("F19").Value = 31.37331497503585
$ws.Range("G19").Value = 3.634734675685971
$ws.Range("J19").Value = 10.99673091871494
$ws.Range("M19").Value = 18.17883222693281
$ws.Range("O19").Value = 22.82058983672136

$ws.Range("B20").Value = 15.97005901335878
$ws.Range("C20").Value = 9.803319250541062
$ws.Range("D20").Value = 10.57312052235754
$ws.Range("F20").Value = 31.35565461160457
$ws.Range("G20").Value = 3.633939210043549
$ws.Range("J20").Value = 10.98597736092695
$ws.Range("M20").Value = 18.26306571454598
$ws.Range("O20").Value = 22.79161888711859

$ws.Range("B21").Value = 16.66826843487963
$ws.Range("C21").Value = 10.34707830166611
$ws.Range("D21").Value = 10.61549267578497
$ws.Range("F21").Value = 31.3082467674195
$ws.Range("G21").Value = 3.631350802093093
$ws.Range("J21").Value = 10.95181851455632
$ws.Range("M21").Value = 18.54572341515569
$ws.Range("O21").Value = 22.70261301524514

$ws.Range("B22").Value = 17.10946507628657
$ws.Range("C22").Value = 10.68647975231574
$ws.Range("D22").Value = 10.64430757940778
$ws.Range("F22").Value = 31.28618860916732
$ws.Range("G22").Value = 3.629721632916341
$ws.Range("J22").Value = 10.93096392063869
$ws.Range("M22").Value = 18.72993678068148
$ws.Range("O22").Value = 22.65069072580643

$ws.Range("B23").Value = 16.87540325401888
$ws.Range("C23").Value = 10.50679762215468
$ws.Range("D23").Value = 10.62882745837187
$ws.Range("F23").Value = 31.29714269161783
$ws.Range("G23").Value = 3.630585450247192
$ws.Range("J23").Value = 10.94195984266049
$ws.Range("M23").Value = 18.63168239301512
$ws.Range("O23").Value = 22.67782853112588

$ws.Range("B24").Value = 15.95839880446796
$ws.Range("C24").Value = 9.794162379641628
$ws.Range("D24").Value = 10.57244727304965
$ws.Range("F24").Value = 31.35657965516546
$ws.Range("G24").Value = 3.633982575541468
$ws.Range("J24").Value = 10.98656047935631
$ws.Range("M24").Value = 18.25844069390066
$ws.Range("O24").Value = 22.79317858387317

$ws.Range("B25").Value = 14.90773702867804
$ws.Range("C25").Value = 8.957070853715898
$ws.Range("D25").Value = 10.51660737738579
$ws.Range("F25").Value = 31.4587718109505
$ws.Range("G25").Value = 3.637917150280602
$ws.Range("J25").Value = 11.04097709192725
$ws.Range("M25").Value = 17.85527086173743
$ws.Range("O25").Value = 22.94417986264858
